# Add a new header row to each PID-list worksheet (WG, MAF, IGN).
# The header documents the meaning of each existing column:
# Name, Unit, Equation, Format, Address, Length, Signed, ProgMin,
# ProgMax, WarnMin, WarnMax, Smoothing, Enabled.

$wb = $excel.ActiveWorkbook

$headers = @("Name","Unit","Equation","Format","Address","Length","Signed","ProgMin","ProgMax","WarnMin","WarnMax","Smoothing","Enabled")

$sheetCount = $wb.Worksheets.Count

for ($s = 1; $s -le $sheetCount; $s++) {
    $ws = $wb.Worksheets.Item($s)

    # Insert a fresh blank row at the top, shifting all existing data down.
    [void]$ws.Rows.Item(1).Insert()

    # Populate the new row with the column headers.
    for ($i = 0; $i -lt $headers.Length; $i++) {
        $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
    }

    # Select the full used range of the sheet, matching the post-edit workbook state.
    $lastRow = $ws.Cells.Item(1, 1).End(-4121).Row
    [void]$ws.Range("A1:M" + $lastRow).Select()
}

# Re-activate the sheet that was active before the edit (IGN, the 3rd tab).
[void]$wb.Worksheets.Item(3).Activate()
[void]$wb.Worksheets.Item(3).Range("A1:M1").Select()
